$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Ind. Rivadavia vs River Plate (odds tweak) ---
$ws.Range('A2').Value = 'b9JYT7LE'
$ws.Range('B2').Value = '21/11/2024'
$ws.Range('C2').Value = '21:30'
$ws.Range('D2').Value = 'ARGENTINA - TORNEO BETANO'
$ws.Range('E2').Value = 'Ind. Rivadavia'
$ws.Range('F2').Value = 'River Plate'
$ws.Range('G2').Value = 8.5
$ws.Range('H2').Value = 4
$ws.Range('I2').Value = 1.45
$ws.Range('J2').Value = 8
$ws.Range('K2').Value = 2.2
$ws.Range('L2').Value = 2.05
$ws.Range('M2').Value = 1.08
$ws.Range('N2').Value = 8
$ws.Range('O2').Value = 1.33
$ws.Range('P2').Value = 3.25
$ws.Range('Q2').Value = 2.1
$ws.Range('R2').Value = 1.7
$ws.Range('S2').Value = 1.44
$ws.Range('T2').Value = 2.63
$ws.Range('U2').Value = 2.25
$ws.Range('V2').Value = 1.57
$ws.Range('W2').Value = 17
$ws.Range('X2').Value = 41
$ws.Range('Y2').Value = 26
$ws.Range('Z2').Value = 101
$ws.Range('AA2').Value = 67
$ws.Range('AB2').Value = 67
$ws.Range('AC2').Value = 8
$ws.Range('AD2').Value = 8
$ws.Range('AE2').Value = 23
$ws.Range('AF2').Value = 81
$ws.Range('AG2').Value = 1250
$ws.Range('AH2').Value = 5.5
$ws.Range('AI2').Value = 6
$ws.Range('AJ2').Value = 9
$ws.Range('AK2').Value = 9
$ws.Range('AL2').Value = 15
$ws.Range('AM2').Value = 34
$ws.Range('AN2').Value = 8.5
$ws.Range('AO2').Value = 41
$ws.Range('AP2').Value = 51
$ws.Range('AQ2').Value = 201
$ws.Range('AR2').Value = 251
$ws.Range('AS2').Value = 2.63
$ws.Range('AT2').Value = 10
$ws.Range('AU2').Value = 81
$ws.Range('AV2').Value = 3.2
$ws.Range('AW2').Value = 7.5
$ws.Range('AX2').Value = 23
$ws.Range('AY2').Value = 23
$ws.Range('AZ2').Value = 51
$ws.Range('BA2').Value = 201
$ws.Range('BB2').Value = 501
$ws.Range('BC2').Value = 126
$ws.Range('BD2').Value = 151

# --- Row 4: Junior vs America De Cali (was Once Caldas vs Deportes Tolima) ---
$ws.Range('A4').Value = 'QuRAYZIt'
$ws.Range('B4').Value = '21/11/2024'
$ws.Range('C4').Value = '22:30'
$ws.Range('D4').Value = 'COLOMBIA - PRIMERA A'
$ws.Range('E4').Value = 'Junior'
$ws.Range('F4').Value = 'America De Cali'
$ws.Range('G4').Value = 2.15
$ws.Range('H4').Value = 3.1
$ws.Range('I4').Value = 3.75
$ws.Range('J4').Value = 3
$ws.Range('K4').Value = 1.95
$ws.Range('L4').Value = 4.33
$ws.Range('M4').Value = 1.1
$ws.Range('N4').Value = 7
$ws.Range('O4').Value = 1.44
$ws.Range('P4').Value = 2.63
$ws.Range('Q4').Value = 2.4
$ws.Range('R4').Value = 1.53
$ws.Range('S4').Value = 1.53
$ws.Range('T4').Value = 2.38
$ws.Range('U4').Value = 2
$ws.Range('V4').Value = 1.73
$ws.Range('W4').Value = 6.5
$ws.Range('X4').Value = 9.5
$ws.Range('Y4').Value = 9.5
$ws.Range('Z4').Value = 21
$ws.Range('AA4').Value = 21
$ws.Range('AB4').Value = 34
$ws.Range('AC4').Value = 7
$ws.Range('AD4').Value = 6
$ws.Range('AE4').Value = 17
$ws.Range('AF4').Value = 67
$ws.Range('AG4').Value = 1250
$ws.Range('AH4').Value = 8.5
$ws.Range('AI4').Value = 17
$ws.Range('AJ4').Value = 13
$ws.Range('AK4').Value = 41
$ws.Range('AL4').Value = 34
$ws.Range('AM4').Value = 41
$ws.Range('AN4').Value = 4
$ws.Range('AO4').Value = 13
$ws.Range('AP4').Value = 26
$ws.Range('AQ4').Value = 41
$ws.Range('AR4').Value = 81
$ws.Range('AS4').Value = 2.38
$ws.Range('AT4').Value = 9
$ws.Range('AU4').Value = 67
$ws.Range('AV4').Value = 5.5
$ws.Range('AW4').Value = 21
$ws.Range('AX4').Value = 34
$ws.Range('AY4').Value = 67
$ws.Range('AZ4').Value = 101
$ws.Range('BA4').Value = 301
$ws.Range('BB4').Value = 251
$ws.Range('BC4').Value = 126
$ws.Range('BD4').Value = 126

# --- Row 5: Guadalajara Chivas vs Atlas (was Junior vs America De Cali) ---
$ws.Range('A5').Value = '6cNu9v1t'
$ws.Range('B5').Value = '21/11/2024'
$ws.Range('C5').Value = '22:05'
$ws.Range('D5').Value = 'MEXICO - LIGA MX'
$ws.Range('E5').Value = 'Guadalajara Chivas'
$ws.Range('F5').Value = 'Atlas'
$ws.Range('G5').Value = 1.73
$ws.Range('H5').Value = 3.75
$ws.Range('I5').Value = 5
$ws.Range('J5').Value = 2.38
$ws.Range('K5').Value = 2.1
$ws.Range('L5').Value = 5.5
$ws.Range('M5').Value = 1.07
$ws.Range('N5').Value = 9
$ws.Range('O5').Value = 1.4
$ws.Range('P5').Value = 2.75
$ws.Range('Q5').Value = 2.25
$ws.Range('R5').Value = 1.62
$ws.Range('S5').Value = 1.5
$ws.Range('T5').Value = 2.5
$ws.Range('U5').Value = 2.1
$ws.Range('V5').Value = 1.67
$ws.Range('W5').Value = 6
$ws.Range('X5').Value = 7
$ws.Range('Y5').Value = 9
$ws.Range('Z5').Value = 13
$ws.Range('AA5').Value = 15
$ws.Range('AB5').Value = 34
$ws.Range('AC5').Value = 8.5
$ws.Range('AD5').Value = 7
$ws.Range('AE5').Value = 21
$ws.Range('AF5').Value = 67
$ws.Range('AG5').Value = 501
$ws.Range('AH5').Value = 11
$ws.Range('AI5').Value = 23
$ws.Range('AJ5').Value = 17
$ws.Range('AK5').Value = 51
$ws.Range('AL5').Value = 41
$ws.Range('AM5').Value = 51
$ws.Range('AN5').Value = 3.5
$ws.Range('AO5').Value = 9
$ws.Range('AP5').Value = 23
$ws.Range('AQ5').Value = 34
$ws.Range('AR5').Value = 51
$ws.Range('AS5').Value = 2.5
$ws.Range('AT5').Value = 9
$ws.Range('AU5').Value = 67
$ws.Range('AV5').Value = 6.5
$ws.Range('AW5').Value = 29
$ws.Range('AX5').Value = 41
$ws.Range('AY5').Value = 101
$ws.Range('AZ5').Value = 151
$ws.Range('BA5').Value = 351
$ws.Range('BB5').Value = 201
$ws.Range('BC5').Value = 51
$ws.Range('BD5').Value = 51

# --- Remove now-duplicate row 6 (Guadalajara Chivas vs Atlas moved to row 5) ---
$ws.Rows.Item(6).Delete()

